$wb = $excel.ActiveWorkbook

# Update recalculated p-values in column C across histology sheets
# (values reflect removal of derived cell-line-only samples from the analysis)

$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Range("C3").Value = 0.609339066093391
$ws.Range("C6").Value = 0.0625937406259374
$ws.Range("C7").Value = 0.224777522247775
$ws.Range("C8").Value = 0.776322367763224

$ws = $wb.Worksheets.Item("Non-neoplastic tumor")
$ws.Range("C3").Value = 0.906709329067093
$ws.Range("C5").Value = 0.0000999900009999
$ws.Range("C6").Value = 0.245875412458754
$ws.Range("C7").Value = 0.856514348565144

$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Range("C3").Value = 0.188981101889811
$ws.Range("C6").Value = 0.518748125187481
$ws.Range("C7").Value = 0.79002099790021
$ws.Range("C8").Value = 0.6996300369963

$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Range("C3").Value = 0.781921807819218
$ws.Range("C6").Value = 0.468853114688531
$ws.Range("C7").Value = 0.302269773022698
$ws.Range("C8").Value = 0.67973202679732

$ws = $wb.Worksheets.Item("Schwannoma")
$ws.Range("C3").Value = 0.807419258074193
$ws.Range("C5").Value = 0.508049195080492
$ws.Range("C6").Value = 0.932806719328067
$ws.Range("C7").Value = 0.848315168483152

$ws = $wb.Worksheets.Item("Mesenchymal tumor")
$ws.Range("C3").Value = 0.120887911208879
$ws.Range("C5").Value = 0.0003999600039996
$ws.Range("C6").Value = 0.0212978702129787
$ws.Range("C7").Value = 0.842115788421158
$ws.Range("C8").Value = 0.485451454854515

$ws = $wb.Worksheets.Item("Germ cell tumor")
$ws.Range("C3").Value = 0.614438556144386
$ws.Range("C4").Value = 0.0000999900009999
$ws.Range("C5").Value = 0.0108989101089891
$ws.Range("C6").Value = 0.859414058594141
$ws.Range("C7").Value = 0.354564543545645
$ws.Range("C8").Value = 0.279172082791721

$ws = $wb.Worksheets.Item("Craniopharyngioma")
$ws.Range("C3").Value = 0.739126087391261
$ws.Range("C5").Value = 0.0057994200579942
$ws.Range("C6").Value = 0.163483651634837
$ws.Range("C7").Value = 0.576942305769423

$ws = $wb.Worksheets.Item("Other tumor")
$ws.Range("C3").Value = 0.0435956404359564
$ws.Range("C5").Value = 0.0004999500049995
$ws.Range("C6").Value = 0.877812218778122
$ws.Range("C7").Value = 0.993000699930007

$ws = $wb.Worksheets.Item("Ependymoma")
$ws.Range("C3").Value = 0.406659334066593
$ws.Range("C6").Value = 0.264073592640736
$ws.Range("C7").Value = 0.280571942805719
$ws.Range("C8").Value = 0.965103489651035

$ws = $wb.Worksheets.Item("DIPG or DMG")
$ws.Range("C2").Value = 0.579588377492154
$ws.Range("C3").Value = 0.618138186181382
$ws.Range("C6").Value = 0.0143985601439856
$ws.Range("C7").Value = 0.764823517648235
$ws.Range("C8").Value = 0.111213887506894
$ws.Range("C9").Value = 0.239428006470177

$ws = $wb.Worksheets.Item("ATRT")
$ws.Range("C3").Value = 0.367963203679632
$ws.Range("C6").Value = 0.558444155584442
$ws.Range("C7").Value = 0.279172082791721
$ws.Range("C8").Value = 0.172282771722828

$ws = $wb.Worksheets.Item("Other high-grade glioma")
$ws.Range("C2").Value = 0.967518506841397
$ws.Range("C3").Value = 0.585041495850415
$ws.Range("C6").Value = 0.890810918908109
$ws.Range("C7").Value = 0.0096990300969903
$ws.Range("C8").Value = 0.64043595640436
$ws.Range("C9").Value = 0.420853921755544
$ws.Range("C10").Value = 0.976650822576724

$ws = $wb.Worksheets.Item("Meningioma")
$ws.Range("C3").Value = 0.947405259474053
$ws.Range("C6").Value = 0.673132686731327
$ws.Range("C7").Value = 0.329167083291671

$ws = $wb.Worksheets.Item("Neurofibroma plexiform")
$ws.Range("C4").Value = 0.0005999400059994
$ws.Range("C6").Value = 0.0585941405859414
$ws.Range("C7").Value = 0.302069793020698

$ws = $wb.Worksheets.Item("Oligodendroglioma")
$ws.Range("C3").Value = 0.200979902009799
$ws.Range("C4").Value = 0.0001999800019998
$ws.Range("C5").Value = 0.472652734726527
$ws.Range("C7").Value = 0.37986201379862
